# Reorders the "EC" (Estado de Cuenta) detail table on Hoja1 (rows 16-36):
# instead of being grouped by period with workers interleaved, the rows are
# now grouped by worker (in original first-appearance order), and within
# each worker ordered by period descending (2001, 1912, 1911, ..., 1907).
# The underlying (worker, period) -> (Valor Mora, Salario Basico) data is
# unchanged; only row order / cell contents are rewritten in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$rows = @(
    @{ Row = 16; Doc = "64521394";    Nombre = "DIANI BERRIO DE DEL TORO";           Periodo = "2001"; Mora = 21874; Salario = 781242 },
    @{ Row = 17; Doc = "64521394";    Nombre = "DIANI BERRIO DE DEL TORO";           Periodo = "1912"; Mora = 31249; Salario = 781242 },
    @{ Row = 18; Doc = "64521394";    Nombre = "DIANI BERRIO DE DEL TORO";           Periodo = "1911"; Mora = 31249; Salario = 781242 },
    @{ Row = 19; Doc = "64521394";    Nombre = "DIANI BERRIO DE DEL TORO";           Periodo = "1910"; Mora = 31249; Salario = 781242 },
    @{ Row = 20; Doc = "64521394";    Nombre = "DIANI BERRIO DE DEL TORO";           Periodo = "1909"; Mora = 31249; Salario = 781242 },
    @{ Row = 21; Doc = "64521394";    Nombre = "DIANI BERRIO DE DEL TORO";           Periodo = "1908"; Mora = 31249; Salario = 781242 },
    @{ Row = 22; Doc = "64521394";    Nombre = "DIANI BERRIO DE DEL TORO";           Periodo = "1907"; Mora = 33125; Salario = 781242 },
    @{ Row = 23; Doc = "1143381314";  Nombre = "ALFONSO ALVAREZ GARCIA";             Periodo = "2001"; Mora = 21874; Salario = 781242 },
    @{ Row = 24; Doc = "1143381314";  Nombre = "ALFONSO ALVAREZ GARCIA";             Periodo = "1912"; Mora = 31249; Salario = 781242 },
    @{ Row = 25; Doc = "1143381314";  Nombre = "ALFONSO ALVAREZ GARCIA";             Periodo = "1911"; Mora = 33125; Salario = 781242 },
    @{ Row = 26; Doc = "1143381314";  Nombre = "ALFONSO ALVAREZ GARCIA";             Periodo = "1910"; Mora = 33125; Salario = 781242 },
    @{ Row = 27; Doc = "1143381314";  Nombre = "ALFONSO ALVAREZ GARCIA";             Periodo = "1909"; Mora = 33125; Salario = 781242 },
    @{ Row = 28; Doc = "1143381314";  Nombre = "ALFONSO ALVAREZ GARCIA";             Periodo = "1908"; Mora = 33125; Salario = 781242 },
    @{ Row = 29; Doc = "1143381314";  Nombre = "ALFONSO ALVAREZ GARCIA";             Periodo = "1907"; Mora = 33125; Salario = 781242 },
    @{ Row = 30; Doc = "33272637";    Nombre = "SHIRLEY PATRICIA DE ORO PALACIN";    Periodo = "2001"; Mora = 21874; Salario = 828116 },
    @{ Row = 31; Doc = "33272637";    Nombre = "SHIRLEY PATRICIA DE ORO PALACIN";    Periodo = "1912"; Mora = 31249; Salario = 828116 },
    @{ Row = 32; Doc = "33272637";    Nombre = "SHIRLEY PATRICIA DE ORO PALACIN";    Periodo = "1911"; Mora = 33125; Salario = 828116 },
    @{ Row = 33; Doc = "33272637";    Nombre = "SHIRLEY PATRICIA DE ORO PALACIN";    Periodo = "1910"; Mora = 33125; Salario = 828116 },
    @{ Row = 34; Doc = "33272637";    Nombre = "SHIRLEY PATRICIA DE ORO PALACIN";    Periodo = "1909"; Mora = 33125; Salario = 828116 },
    @{ Row = 35; Doc = "33272637";    Nombre = "SHIRLEY PATRICIA DE ORO PALACIN";    Periodo = "1908"; Mora = 33125; Salario = 828116 },
    @{ Row = 36; Doc = "33272637";    Nombre = "SHIRLEY PATRICIA DE ORO PALACIN";    Periodo = "1907"; Mora = 33125; Salario = 828116 }
)

foreach ($item in $rows) {
    $r = $item.Row
    $ws.Cells.Item($r, 3).Value = $item.Doc
    $ws.Cells.Item($r, 4).Value = $item.Nombre
    $ws.Cells.Item($r, 5).Value = $item.Periodo
    $ws.Cells.Item($r, 6).Value = $item.Mora
    $ws.Cells.Item($r, 7).Value = $item.Salario
}
